$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "3 / 4 V's [Volume, Velocity, Variety, Veracity]"
#       -> "3 V's [Volume, Velocity, Variety]"
#    (the " / 4" is dropped and the trailing ", Veracity" is dropped)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("3 / 4 V", $true, $false, $false, $false, $false, $true, 1, $false, "3 V", 2) | Out-Null

$d.Content.Find.Execute("Variety, Veracity]", $true, $false, $false, $false, $false, $true, 1, $false, "Variety]", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Future for Big Data " + [proofErr]etc[proofErr]
#       -> "Future for Big Data etc" as a single clean run (no proofErr marks)
#    Rebuild the paragraph content from scratch so the stale spell-check
#    markers are dropped entirely.
# ---------------------------------------------------------------------------
$rngFuture = $d.Content
$rngFuture.Find.Execute("Future for Big Data", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$futurePara = $rngFuture.Paragraphs(1)
$paraStart = $futurePara.Range.Start
$paraEnd = $futurePara.Range.End

# Insert a fresh paragraph (inherits the same pPr/numbering) right before it ...
$d.Range($paraStart, $paraStart).InsertParagraphBefore() | Out-Null

# ... fill it with the clean replacement text ...
$newText = "Future for Big Data etc"
$d.Range($paraStart, $paraStart).Text = $newText

# ... then delete the old paragraph (including its proofErr markers) which has
# now been shifted along by the new paragraph mark + new text length.
$shift = 1 + $newText.Length
$d.Range($paraStart + $shift, $paraEnd + $shift).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) Conclusion paragraph: the visible text doesn't change, but it ends up
#    split into three runs instead of two:
#      "Here, I will off" | "er ... summarising wha" | "t's been covered ..."
#    Use (temporary) bookmarks to force the run split at those two positions
#    without altering any text, then remove the temporary bookmarks again.
# ---------------------------------------------------------------------------
$rngSplit1 = $d.Content
$rngSplit1.Find.Execute("Here, I will off", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngSplit1.Collapse(0) | Out-Null
$d.Bookmarks.Add("TempSplitMark1", $rngSplit1) | Out-Null

$rngSplit2 = $d.Content
$rngSplit2.Find.Execute("summarising wha", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngSplit2.Collapse(0) | Out-Null
$d.Bookmarks.Add("TempSplitMark2", $rngSplit2) | Out-Null

$d.Bookmarks("TempSplitMark1").Delete()
$d.Bookmarks("TempSplitMark2").Delete()

# ---------------------------------------------------------------------------
# 4) Move the (hidden) "_GoBack" last-edit-position bookmark from the end of
#    the document to right after "Variety" in the first paragraph we edited
#    above (i.e. right before the closing "]").
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$rngGoBack = $d.Content
$rngGoBack.Find.Execute("Variety", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngGoBack.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $rngGoBack) | Out-Null
